$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the top parameter values (hipLen, toeLen) ---
$ws.Range("B1").Value = 2.0754000000000001
$ws.Range("B2").Value = 2.5745

# --- Remove the human_9_exo model column (old column D) ---
# This shifts human_9_load_exo (old column E) left into column D,
# and shrinks the used range from AE to AD.
$ws.Columns("D").Delete()

# --- New simulation run timestamps for the "forward"/"backward" blocks ---
# (written in the same cell order the original log used, so the shared
# string table comes out in the same sequence as the authored workbook)
$ws.Range("B15").Value = "0206000141"
$ws.Range("B9").Value  = "0205235730"
$ws.Range("B14").Value = "0206000110"
$ws.Range("B13").Value = "0206000054"
$ws.Range("B11").Value = "0206000033"
$ws.Range("B10").Value = "0207222945"
$ws.Range("B21").Value = "0207223949"
$ws.Range("B19").Value = "0207223824"
$ws.Range("B23").Value = "0207223659"
$ws.Range("B20").Value = "0207232750"
$ws.Range("B12").Value = "0207233003"
$ws.Range("B24").Value = "0207233150"
$ws.Range("B25").Value = "0207233103"
$ws.Range("B22").Value = "0207234726"

$ws.Range("C14").Value = "0207235429"
$ws.Range("C9").Value  = "0207235223"
$ws.Range("C12").Value = "0207235341"
$ws.Range("C15").Value = "0207235354"
$ws.Range("C13").Value = "0207235425"
$ws.Range("C11").Value = "0207235328"
$ws.Range("C10").Value = "0207235314"
$ws.Range("C20").Value = "0208004556"
$ws.Range("C21").Value = "0208004533"
$ws.Range("C19").Value = "0208004350"
$ws.Range("C23").Value = "0208004603"
$ws.Range("C25").Value = "0208004618"
$ws.Range("C22").Value = "0208004310"
$ws.Range("C24").Value = "0208004427"

$ws.Range("D13").Value = "0208011030"
$ws.Range("D10").Value = "0208011059"
$ws.Range("D11").Value = "0208010942"
$ws.Range("D15").Value = "0208010951"
$ws.Range("D12").Value = "0208010828"
$ws.Range("D14").Value = "0208110847"
$ws.Range("D9").Value  = "0208110754"
$ws.Range("D22").Value = "0208111153"
$ws.Range("D23").Value = "0208110658"
$ws.Range("D25").Value = "0208111009"
$ws.Range("D20").Value = "0208110907"
$ws.Range("D21").Value = "0208112819"
$ws.Range("D19").Value = "0208113044"
$ws.Range("D24").Value = "0208112646"

# --- Move the active selection to D24 (matches where the last edit landed) ---
$ws.Range("D24").Select()
